$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 value update
$ws.Range("B3").Value = 0.9215856593215532

# Rows 4 and 5: swap labels (x <-> y) and update values
$ws.Range("A4").Value = "y"
$ws.Range("B4").Value = 0.8888175447215017

$ws.Range("A5").Value = "x"
$ws.Range("B5").Value = 0.8872212110854774

# Remaining B column value updates (labels unchanged)
$ws.Range("B6").Value = 0.8821066056491267
$ws.Range("B7").Value = 0.1282727869737869
$ws.Range("B8").Value = 0.1272413322081764
$ws.Range("B9").Value = 0.09719151518199726
$ws.Range("B10").Value = 0.09525504563025693
$ws.Range("B11").Value = 0.08188540924334869
$ws.Range("B12").Value = 0.05873119343477257
$ws.Range("B13").Value = 0.01876460707401062
$ws.Range("B14").Value = 0.008963348666493816
$ws.Range("B15").Value = 0.008511647393192818
$ws.Range("B16").Value = 0.006862420000157086
$ws.Range("B17").Value = -0.0001394779193105116
$ws.Range("B18").Value = -0.0003574708906823419
$ws.Range("B19").Value = -0.001094085972133455
$ws.Range("B20").Value = -0.009677735194627292
$ws.Range("B21").Value = -0.0106525476532575
$ws.Range("B22").Value = -0.02393224913792053
$ws.Range("B23").Value = -0.04954682179790952
$ws.Range("B24").Value = -0.0522698070334784
$ws.Range("B25").Value = -0.07252059251127617
$ws.Range("B26").Value = -0.09544297924939737
$ws.Range("B27").Value = -0.0969723330799748
$ws.Range("B28").Value = -0.1009108224223475
$ws.Range("B29").Value = -0.3070850070127243
